$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap the header values of B1 ("Longitude") and C1 ("Latitide")
$b1 = $ws.Range("B1").Value2
$c1 = $ws.Range("C1").Value2
$ws.Range("B1").Value2 = $c1
$ws.Range("C1").Value2 = $b1

# Move the active cell selection from B5 to D5
$ws.Range("D5").Select()
